# repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect the repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    7  = 4
    8  = -4
    12 = 0
    15 = -4
    18 = -1
    24 = -1
    25 = 3
    32 = 2
    37 = 0
    39 = 0
    40 = 2
    42 = 1
    44 = 0
    47 = -4
    52 = 0
    55 = 0
    57 = -4
    61 = -1
    68 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
